# Sheet1 holds the "Inputs" -> "Outputs" demo table used by the wrapper
# test fixture. Fill in the three input cells (x, b, s) that were left
# blank, which also drives a recalculation of the dependent formulas in
# column F (y = x*2.1, bout = IF(b, FALSE, TRUE), sout = LOWER(s)).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# x (B2) -> numeric input, flows into F2 (=x*2.1)
$ws.Range("B2").Value = 12

# b (B3) -> boolean input, flows into F3 (=IF(b, FALSE, TRUE))
$ws.Range("B3").Value = $true

# s (B4) -> string input, flows into F4 (=LOWER(s))
$ws.Range("B4").Value = "Hello"

# Leave the selection on F3, matching the saved view state.
$ws.Range("F3").Select() | Out-Null
